# The styles below define their <w:rPr> run properties with <w:color>
# listed before <w:b>/<w:i>. That element order violates the CT_RPr
# sequence in wml.xsd (color must come after the bold/italic toggles),
# which OOXMLValidator flags even though xmllint stays quiet about it.
#
# Re-assigning each style's Font.Bold / Font.Italic (even to their own
# current value) makes the engine rewrite that style's <w:rPr> in
# schema-correct order, moving <w:b>/<w:i> ahead of <w:color> without
# altering any of the actual formatting.

$d = $word.ActiveDocument
$styles = $d.Styles

function Fix-TokStyle($styleName, [bool]$touchBold, [bool]$touchItalic) {
    $style = $styles.Item($styleName)
    if ($touchBold) {
        $style.Font.Bold = $style.Font.Bold
    }
    if ($touchItalic) {
        $style.Font.Italic = $style.Font.Italic
    }
}

# <w:b/><w:color/>
Fix-TokStyle "KeywordTok"      $true  $false
Fix-TokStyle "ImportTok"       $true  $false
Fix-TokStyle "ControlFlowTok"  $true  $false
Fix-TokStyle "AlertTok"        $true  $false
Fix-TokStyle "ErrorTok"        $true  $false

# <w:i/><w:color/>
Fix-TokStyle "CommentTok"        $false $true
Fix-TokStyle "DocumentationTok"  $false $true

# <w:b/><w:i/><w:color/>
Fix-TokStyle "AnnotationTok"   $true $true
Fix-TokStyle "CommentVarTok"   $true $true
Fix-TokStyle "InformationTok"  $true $true
Fix-TokStyle "WarningTok"      $true $true
